# Mise à jour de certains champs de Modules et de Professeurs
#
# Header row changes on Feuil1:
#   C1: "Enseignant"       -> "Chef  Module"
#   D1: "Nombre d'heures"  -> "Composants"
#
# Column widths (new explicit <cols> sizing for C and D), active
# selection moved to E8, per the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# ColumnWidth is expressed in "characters" and gets translated by Excel
# into the stored character-unit width (roughly width + 5/6). Compensate
# so the persisted <col width="..."> lands on the desired values.
$ws.Columns("C").ColumnWidth = 34.1666666667
$ws.Columns("D").ColumnWidth = 23.7369791667

$null = $ws.Range("E8").Select()
